$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column headers: G = "7 Robinwood" (centered, like the B/C/D headers),
# J = "Sample House" (same shared string as H1/I1, no special style)
$ws.Range("G1").Value = "7 Robinwood"
$ws.Range("G1").HorizontalAlignment = -4108  # xlCenter, matches B1/C1/D1 style
$ws.Range("J1").Value = "Sample House"

# Per-year data for the new "7 Robinwood" (G) and "Sample House" (J) columns
$gvals = @(13058, 13358, 13658, 13958, 14258, 14558, 14858, 15158, 15458, 15758)
$jvals = @(15850, 16150, 16450, 16750, 17050, 17350, 17650, 17950, 18250, 18550)

for ($i = 0; $i -lt 10; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $gvals[$i]
    $ws.Cells.Item($row, 10).Value = $jvals[$i]
}

# Averages row
$ws.Range("G13").Formula = "=AVERAGE(G2:G12)"
$ws.Range("J13").Formula = "=AVERAGE(J2:J12)"

# Best-fit-ish width for column I (bystander column whose width gets
# recorded once the new columns are auto-fit alongside it)
$ws.Columns.Item(9).ColumnWidth = 11.998697916666666

# Zoom level changed by the author while reviewing the new columns
$ws.Application.ActiveWindow.Zoom = 139
